$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new "through" date
$ws.Name = "Through 2022-08-16"

# Update the header label cell (A1's paired total header, "2022 (through 08-15)")
$ws.Range("I1").Value = "2022 (through 08-16)"

# Update the September (row 9) total for 2022 column
$ws.Range("I9").Value = 90

# Update the grand Total row (row 14) for 2022 column
$ws.Range("I14").Value = 1061
